$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 269
$ws.Range("H2").Value = "kitchens"
$ws.Range("I2").Value = "target"
$ws.Range("K2").Value = "j"
$ws.Range("L2").Value = "stimuli/img_yeh72.png"
$ws.Range("M2").Value = 68.66666666666667
$ws.Range("N2").Value = 45.21212121212121
$ws.Range("O2").Value = 56.93939393939394
$ws.Range("P2").Value = 33
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = 4

$ws.Range("F3").Value = 270
$ws.Range("H3").Value = "kitchens"
$ws.Range("I3").Value = "target"
$ws.Range("K3").Value = "j"
$ws.Range("L3").Value = "stimuli/img_wyl6z.png"
$ws.Range("M3").Value = 59.8235294117647
$ws.Range("N3").Value = 36.23529411764706
$ws.Range("O3").Value = 48.02941176470588
$ws.Range("P3").Value = 34
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = 3
$ws.Range("S3").Value = 3

$ws.Range("F4").Value = 271
$ws.Range("H4").Value = "kitchens"
$ws.Range("I4").Value = "target"
$ws.Range("K4").Value = "j"
$ws.Range("L4").Value = "stimuli/img_es7o2.png"
$ws.Range("M4").Value = 52.48571428571429
$ws.Range("N4").Value = 27.54285714285714
$ws.Range("O4").Value = 40.01428571428572
$ws.Range("P4").Value = 35
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 2

$ws.Range("F5").Value = 272
$ws.Range("H5").Value = "living_rooms"
$ws.Range("I5").Value = "distractor"
$ws.Range("K5").Value = "f"
$ws.Range("L5").Value = "stimuli/img_37hgm.png"
$ws.Range("M5").Value = 70.95454545454545
$ws.Range("N5").Value = 54.77272727272727
$ws.Range("O5").Value = 62.86363636363636
$ws.Range("P5").Value = 44
$ws.Range("Q5").Value = 6
$ws.Range("R5").Value = 6
$ws.Range("S5").Value = 6

$ws.Range("F6").Value = 273
$ws.Range("H6").Value = "kitchens"
$ws.Range("I6").Value = "target"
$ws.Range("K6").Value = "j"
$ws.Range("L6").Value = "stimuli/img_d8xbu.png"
$ws.Range("M6").Value = 91.36363636363636
$ws.Range("N6").Value = 73.18181818181819
$ws.Range("O6").Value = 82.27272727272728
$ws.Range("P6").Value = 33
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = 10
$ws.Range("S6").Value = 10

$ws.Range("F7").Value = 274
$ws.Range("H7").Value = "kitchens"
$ws.Range("I7").Value = "target"
$ws.Range("K7").Value = "j"
$ws.Range("L7").Value = "stimuli/img_eatdk.png"
$ws.Range("M7").Value = 81.40625
$ws.Range("N7").Value = 61.375
$ws.Range("O7").Value = 71.390625
$ws.Range("P7").Value = 32
$ws.Range("Q7").Value = 8
$ws.Range("R7").Value = 8
$ws.Range("S7").Value = 8

$ws.Range("F8").Value = 275
$ws.Range("H8").Value = "kitchens"
$ws.Range("I8").Value = "target"
$ws.Range("K8").Value = "j"
$ws.Range("L8").Value = "stimuli/img_9mky8.png"
$ws.Range("M8").Value = 84.32352941176471
$ws.Range("N8").Value = 65.17647058823529
$ws.Range("O8").Value = 74.75
$ws.Range("P8").Value = 34
$ws.Range("Q8").Value = 9
$ws.Range("R8").Value = 9
$ws.Range("S8").Value = 9

$ws.Range("F9").Value = 276
$ws.Range("H9").Value = "bedrooms"
$ws.Range("I9").Value = "distractor"
$ws.Range("K9").Value = "f"
$ws.Range("L9").Value = "stimuli/img_n9xll.png"
$ws.Range("M9").Value = 77.14285714285714
$ws.Range("N9").Value = 59.21428571428572
$ws.Range("O9").Value = 68.17857142857143
$ws.Range("P9").Value = 42
$ws.Range("Q9").Value = 7
$ws.Range("R9").Value = 7
$ws.Range("S9").Value = 7

$ws.Range("F10").Value = 277
$ws.Range("H10").Value = "kitchens"
$ws.Range("I10").Value = "target"
$ws.Range("K10").Value = "j"
$ws.Range("L10").Value = "stimuli/img_iyxnj.png"
$ws.Range("M10").Value = 75.30555555555556
$ws.Range("N10").Value = 54.33333333333334
$ws.Range("O10").Value = 64.81944444444444
$ws.Range("P10").Value = 36
$ws.Range("Q10").Value = 6
$ws.Range("R10").Value = 6
$ws.Range("S10").Value = 6

$ws.Range("F11").Value = 278
$ws.Range("H11").Value = "kitchens"
$ws.Range("I11").Value = "target"
$ws.Range("K11").Value = "j"
$ws.Range("L11").Value = "stimuli/img_q1ynd.png"
$ws.Range("M11").Value = 70.05714285714286
$ws.Range("N11").Value = 47.31428571428572
$ws.Range("O11").Value = 58.68571428571429
$ws.Range("P11").Value = 35
$ws.Range("Q11").Value = 5
$ws.Range("R11").Value = 5
$ws.Range("S11").Value = 5

$ws.Range("F12").Value = 279
$ws.Range("H12").Value = "kitchens"
$ws.Range("I12").Value = "target"
$ws.Range("K12").Value = "j"
$ws.Range("L12").Value = "stimuli/img_ce9vx.png"
$ws.Range("M12").Value = 75.90909090909091
$ws.Range("N12").Value = 57.12121212121212
$ws.Range("O12").Value = 66.51515151515152
$ws.Range("P12").Value = 33
$ws.Range("Q12").Value = 7
$ws.Range("R12").Value = 7
$ws.Range("S12").Value = 7

$ws.Range("F13").Value = 280
$ws.Range("H13").Value = "kitchens"
$ws.Range("I13").Value = "target"
$ws.Range("K13").Value = "j"
$ws.Range("L13").Value = "stimuli/img_inqod.png"
$ws.Range("M13").Value = 70.84848484848484
$ws.Range("N13").Value = 50.63636363636363
$ws.Range("O13").Value = 60.74242424242424
$ws.Range("P13").Value = 33
$ws.Range("Q13").Value = 5
$ws.Range("R13").Value = 5
$ws.Range("S13").Value = 5

$ws.Range("F14").Value = 281
$ws.Range("H14").Value = "kitchens"
$ws.Range("I14").Value = "target"
$ws.Range("K14").Value = "j"
$ws.Range("L14").Value = "stimuli/img_a8wvq.png"
$ws.Range("M14").Value = 86.25925925925925
$ws.Range("N14").Value = 66.25925925925925
$ws.Range("O14").Value = 76.25925925925925
$ws.Range("P14").Value = 27
$ws.Range("Q14").Value = 10
$ws.Range("R14").Value = 10
$ws.Range("S14").Value = 10

$ws.Range("F15").Value = 282
$ws.Range("H15").Value = "living_rooms"
$ws.Range("I15").Value = "distractor"
$ws.Range("K15").Value = "f"
$ws.Range("L15").Value = "stimuli/img_3sw8t.png"
$ws.Range("M15").Value = 67.48888888888889
$ws.Range("N15").Value = 48.51111111111111
$ws.Range("O15").Value = 58
$ws.Range("P15").Value = 45
$ws.Range("Q15").Value = 5
$ws.Range("R15").Value = 5
$ws.Range("S15").Value = 5

$ws.Range("F16").Value = 283
$ws.Range("H16").Value = "kitchens"
$ws.Range("I16").Value = "target"
$ws.Range("K16").Value = "j"
$ws.Range("L16").Value = "stimuli/img_7wul8.png"
$ws.Range("M16").Value = 43.03030303030303
$ws.Range("N16").Value = 25.54545454545455
$ws.Range("O16").Value = 34.28787878787879
$ws.Range("P16").Value = 33
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = 1
$ws.Range("S16").Value = 1

$ws.Range("F17").Value = 284
$ws.Range("H17").Value = "kitchens"
$ws.Range("I17").Value = "target"
$ws.Range("K17").Value = "j"
$ws.Range("L17").Value = "stimuli/img_60242.png"
$ws.Range("M17").Value = 78.33333333333333
$ws.Range("N17").Value = 57.57575757575758
$ws.Range("O17").Value = 67.95454545454545
$ws.Range("P17").Value = 33
$ws.Range("Q17").Value = 7
$ws.Range("R17").Value = 7
$ws.Range("S17").Value = 7

$ws.Range("F18").Value = 285
$ws.Range("H18").Value = "kitchens"
$ws.Range("I18").Value = "target"
$ws.Range("K18").Value = "j"
$ws.Range("L18").Value = "stimuli/img_79b5l.png"
$ws.Range("M18").Value = 72.74285714285715
$ws.Range("N18").Value = 53.31428571428572
$ws.Range("O18").Value = 63.02857142857143
$ws.Range("P18").Value = 35
$ws.Range("Q18").Value = 6
$ws.Range("R18").Value = 6
$ws.Range("S18").Value = 6

$ws.Range("F19").Value = 286
$ws.Range("H19").Value = "living_rooms"
$ws.Range("I19").Value = "distractor"
$ws.Range("K19").Value = "f"
$ws.Range("L19").Value = "stimuli/img_kq9s9.png"
$ws.Range("M19").Value = 62.30232558139535
$ws.Range("N19").Value = 39.97674418604651
$ws.Range("O19").Value = 51.13953488372093
$ws.Range("P19").Value = 43
$ws.Range("Q19").Value = 4
$ws.Range("R19").Value = 4
$ws.Range("S19").Value = 4

$ws.Range("F20").Value = 287
$ws.Range("H20").Value = "living_rooms"
$ws.Range("I20").Value = "distractor"
$ws.Range("K20").Value = "f"
$ws.Range("L20").Value = "stimuli/img_hmmra.png"
$ws.Range("M20").Value = 54.65853658536585
$ws.Range("N20").Value = 34.24390243902439
$ws.Range("O20").Value = 44.45121951219512
$ws.Range("P20").Value = 41
$ws.Range("Q20").Value = 3
$ws.Range("R20").Value = 3
$ws.Range("S20").Value = 3

$ws.Range("F21").Value = 288
$ws.Range("H21").Value = "kitchens"
$ws.Range("I21").Value = "target"
$ws.Range("K21").Value = "j"
$ws.Range("L21").Value = "stimuli/img_nyv2b.png"
$ws.Range("M21").Value = 11.91176470588235
$ws.Range("N21").Value = 6.852941176470588
$ws.Range("O21").Value = 9.382352941176471
$ws.Range("P21").Value = 34
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = 1
$ws.Range("S21").Value = 1

$ws.Range("F22").Value = 289
$ws.Range("H22").Value = "kitchens"
$ws.Range("I22").Value = "target"
$ws.Range("K22").Value = "j"
$ws.Range("L22").Value = "stimuli/img_ye5sl.png"
$ws.Range("M22").Value = 53.2258064516129
$ws.Range("N22").Value = 34.45161290322581
$ws.Range("O22").Value = 43.83870967741936
$ws.Range("P22").Value = 31
$ws.Range("Q22").Value = 2
$ws.Range("R22").Value = 2
$ws.Range("S22").Value = 2

$ws.Range("F23").Value = 290
$ws.Range("H23").Value = "living_rooms"
$ws.Range("I23").Value = "distractor"
$ws.Range("K23").Value = "f"
$ws.Range("L23").Value = "stimuli/img_iudc4.png"
$ws.Range("M23").Value = 73.625
$ws.Range("N23").Value = 52.275
$ws.Range("O23").Value = 62.95
$ws.Range("P23").Value = 40
$ws.Range("Q23").Value = 6
$ws.Range("R23").Value = 6
$ws.Range("S23").Value = 6

$ws.Range("F24").Value = 291
$ws.Range("H24").Value = "kitchens"
$ws.Range("I24").Value = "target"
$ws.Range("K24").Value = "j"
$ws.Range("L24").Value = "stimuli/img_uwv6y.png"
$ws.Range("M24").Value = 78.88888888888889
$ws.Range("N24").Value = 59.30555555555556
$ws.Range("O24").Value = 69.09722222222223
$ws.Range("P24").Value = 36
$ws.Range("Q24").Value = 8
$ws.Range("R24").Value = 8
$ws.Range("S24").Value = 8

$ws.Range("F25").Value = 292
$ws.Range("H25").Value = "kitchens"
$ws.Range("I25").Value = "target"
$ws.Range("K25").Value = "j"
$ws.Range("L25").Value = "stimuli/img_cv6mf.png"
$ws.Range("M25").Value = 66.8
$ws.Range("N25").Value = 42.08
$ws.Range("O25").Value = 54.44
$ws.Range("P25").Value = 25
$ws.Range("Q25").Value = 4
$ws.Range("R25").Value = 4
$ws.Range("S25").Value = 4

$ws.Range("F26").Value = 293
$ws.Range("H26").Value = "kitchens"
$ws.Range("I26").Value = "target"
$ws.Range("K26").Value = "j"
$ws.Range("L26").Value = "stimuli/img_t90e2.png"
$ws.Range("M26").Value = 83.0625
$ws.Range("N26").Value = 61.96875
$ws.Range("O26").Value = 72.515625
$ws.Range("P26").Value = 32
$ws.Range("Q26").Value = 9
$ws.Range("R26").Value = 9
$ws.Range("S26").Value = 9

$ws.Range("F27").Value = 294
$ws.Range("H27").Value = "kitchens"
$ws.Range("I27").Value = "target"
$ws.Range("K27").Value = "j"
$ws.Range("L27").Value = "stimuli/img_aplao.png"
$ws.Range("M27").Value = 64.09090909090909
$ws.Range("N27").Value = 40.75757575757576
$ws.Range("O27").Value = 52.42424242424242
$ws.Range("P27").Value = 33
$ws.Range("Q27").Value = 3
$ws.Range("R27").Value = 3
$ws.Range("S27").Value = 3
